# v0.6 - encoder failure detection
# Adds a new "trip counter" register (MB_ENCODER_FAIL_TRIPS,) in the 300-series
# block, and two new registers at the end of the 9000-series config block:
# MB_HEARTBEAT_TIMEOUT, and MB_ENCODER_FAIL_TIMEOUT,

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 65 (the pre-existing gap before the 9000-range block),
# which pushes the existing 9000-series rows (previously 66-73) down to 67-74.
$ws.Rows.Item(65).Insert()

$ws.Range("A65").Value = 311
$ws.Range("B65").Value = "MB_ENCODER_FAIL_TRIPS,"
$ws.Range("C65").Value = "Number of trips caused by encoder failure detection"
$ws.Range("D65").Value = "R"

# Append two brand new registers after the existing last row (now row 74).
$ws.Range("A75").Value = 9008
$ws.Range("B75").Value = "MB_HEARTBEAT_TIMEOUT,"
$ws.Range("C75").Value = "Seconds before heartbeat timer expires"
$ws.Range("D75").Value = "R/W"
$ws.Range("E75").Value = "s"

$ws.Range("A76").Value = 9009
$ws.Range("B76").Value = "MB_ENCODER_FAIL_TIMEOUT,"
$ws.Range("C76").Value = "Max milliseconds between encoder pulses before timeout"
$ws.Range("D76").Value = "R/W"
$ws.Range("E76").Value = "ms"

# Match the author's final cursor/selection position.
$ws.Range("E77").Select() | Out-Null
